$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New day block: rows 49-57 (day 41926 / Oct 2014), mirroring the existing
# rows 40-48 block (day 41925).
# ---------------------------------------------------------------------------

# Row 49: date header row
$ws.Range("B49").Value = 41926
$ws.Range("B49").NumberFormat = "mm-dd-yy"
$ws.Range("B49").Interior.Color = 15773696
$ws.Range("C49").Interior.Color = 15773696

# Row 50: Ke hoach (plan)
$ws.Range("B50").Value = "Kế hoạch"
$ws.Range("B50").Interior.Color = 49407
$ws.Range("C50").WrapText = $true
$ws.Range("C50").Value = '`- Công việc 1:Buid giao diện màn hình Công văn.
`- Công việc 2: Kết nối API cho màn hình Công văn'

# Row 51: Ket qua dat duoc (results)
$ws.Range("B51").Value = "Kết quả đạt được"
$ws.Range("B51").Interior.Color = 49407
$ws.Range("C51").WrapText = $true
$ws.Range("C51").Value = '''- Công việc 1: 90%
''- Công việc 2: 100%'

# Row 52: Trang thai (status)
$ws.Range("B52").Value = "Trạng thái"
$ws.Range("B52").Interior.Color = 49407
$ws.Range("C52").Value = 'Chậm tiến độ'

# Row 53: Van de gap phai (problems)
$ws.Range("B53").Value = "Vấn đề gặp phải"
$ws.Range("B53").Interior.Color = 49407
$ws.Range("C53").Value = '''- Vấn đề 1: Chưa crack được photoshop.'

# Row 54: - Van de 2:
$ws.Range("B54").Interior.Color = 49407
$ws.Range("C54").Value = '- Vấn đề 2:'

# Row 55: Giai quyet van de (solutions)
$ws.Range("B55").Value = "Giải quyết vấn đề:"
$ws.Range("B55").Interior.Color = 49407
$ws.Range("C55").Value = '''- Vấn đề 1: Crack photoshop và tiếp tục làm'

# Row 56: - Van de 2: giai quyet nhu sau....
$ws.Range("B56").Interior.Color = 49407
$ws.Range("C56").Value = '- Vấn đề 2: giải quyết như sau....'

# Row 57: Ke hoach ngay mai (tomorrow's plan)
$ws.Range("B57").Value = "Kế hoạch ngày mai"
$ws.Range("B57").Interior.Color = 49407
$ws.Range("C57").WrapText = $true
$ws.Range("C57").Value = 'Hoàn thành màn hình Công văn,
Tiếp tục Buid giao diện màn hình Tìm kiếm.'

# ---------------------------------------------------------------------------
# View state: selection moved to C58, scrolled so row 44 is the top row.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 44
$win.ScrollColumn = 1
$ws.Range("C58").Select()
